# input.xlsx ---> test data for invalid login TC
# InvalidLogin.java --> Script for invalid login TC

$wb = $excel.ActiveWorkbook

# Reset selection on the existing ValidLogin sheet to A1:B2 (was A7:XFD7)
$wsValid = $wb.Worksheets.Item("ValidLogin")
[void]$wsValid.Range("A1:B2").Select()

# Add a new worksheet for invalid-login test data, placed right after ValidLogin
$wsInvalid = $wb.Worksheets.Add($null, $wsValid)
$wsInvalid.Name = "InvalidLogin"

# Populate the InvalidLogin sheet with header row + invalid credentials row
$wsInvalid.Range("A1").Value = "Username"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Select B3 on the new sheet and make it the active/visible tab
[void]$wsInvalid.Range("B3").Select()
$wsInvalid.Activate()

[void]$wb.Save()
